$wb = $excel.ActiveWorkbook
$wsJ = $wb.Worksheets.Item("Journal")
$wsT = $wb.Worksheets.Item("Totaux")

# --- Journal sheet: append two new rows, copying formatting from row 19 ---
$wsJ.Range("A19:E19").Copy($wsJ.Range("A20:E20"))
$wsJ.Range("A19:E19").Copy($wsJ.Range("A21:E21"))

$wsJ.Range("A20").Value = 44978
$wsJ.Range("B20").Value = 3
$wsJ.Range("C20").Value = 0.041666666666666664
$wsJ.Range("D20").Value = "Documentation"
$wsJ.Range("E20").Value = "Analyse des fermetures de ticket dans la partie ""gestion des tckets"""

$wsJ.Range("A21").Value = 44978
$wsJ.Range("B21").Value = 3
$wsJ.Range("C21").Value = 0.020833333333333332
$wsJ.Range("D21").Value = "Analyse"
$wsJ.Range("E21").Value = "Rédaction des problèmes techniques dans le rapport de projet"

$tblJ = $wsJ.ListObjects.Item(1)
$tblJ.Resize($wsJ.Range("A1:E21"))

# --- Totaux sheet: insert a new row before the totals row ---
$wsT.Rows.Item(8).Insert()

$wsT.Range("A8").Value = 44978
$wsT.Range("B8").Formula = "=SUM(Journal!C20:C21)"

$tblT = $wsT.ListObjects.Item(1)
$tblT.Resize($wsT.Range("A1:B9"))

$wsT.Range("B9").Formula = "=SUM(B2:B8)"

# --- Selection / active sheet bookkeeping ---
$wsJ.Range("E21").Select()
$wsT.Activate()
$wsT.Range("E7").Select()
